$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab (sheet name only, file name stays the same)
$ws.Name = "SA"

# Add new row 16 to the averaged-intensities table, matching the formatting
# of row 15 (bold, centered, bordered index cell in column A; shared text
# label in column B; plain numeric values in columns C:P).
$ws.Range("A15:B15").Copy()
$ws.Range("A16:B16").PasteSpecial(-4122)

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C16").Value = 1.085862455819836
$ws.Range("D16").Value = 0.8713262723555942
$ws.Range("E16").Value = 1.024422867435443
$ws.Range("F16").Value = 0.9630349633590677
$ws.Range("G16").Value = 1.085862455819836
$ws.Range("H16").Value = 0.8713262723555942
$ws.Range("I16").Value = 1.040908964541113
$ws.Range("J16").Value = 0.9615677810982491
$ws.Range("K16").Value = 1.029507205846422
$ws.Range("L16").Value = 0.9138907047003079
$ws.Range("M16").Value = 1.085862455819836
$ws.Range("N16").Value = 0.9478745698955184
$ws.Range("O16").Value = 0.9861616397424852
$ws.Range("P16").Value = 0.9863151518945041
